# Update "想去人数" (want-to-go count) figures in the F column on the
# "展览" and "全部类型" sheets, reflecting newer scrape numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 730
$ws1.Range("F4").Value = 246
$ws1.Range("F5").Value = 2747
$ws1.Range("F7").Value = 3717
$ws1.Range("F9").Value = 940

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 730
$ws4.Range("F5").Value = 246
$ws4.Range("F6").Value = 2747
$ws4.Range("F8").Value = 3717
$ws4.Range("F10").Value = 940
